$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.643.02"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -6.75%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.433.87"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -9.27%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "534.01"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.12"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -8.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.568"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.94%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.444.13"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -9.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0983"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -7.09%  "
$ws.Range("E11").Value = "  -2.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.31"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.88%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.346"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -5.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.873.11"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -9.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.53"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -11.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "58.552.68"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -6.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000136"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -6.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.490.36"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -7.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.02"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -6.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.31"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -6.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "321.59"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -6.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.964"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.66"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -8.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.15"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -4.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.447"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -12.72%  "
$ws.Range("E26").Value = "  -6.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.973"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.62"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -6.67%  "
$ws.Range("E29").Value = "  -6.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0762"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -10.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.59"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -9.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.21"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -12.08%  "
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "154.88"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -5.95%  "
$ws.Range("E35").Value = "  -5.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.37"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -6.79%  "
$ws.Range("E37").Value = "  -10.75%  "
$ws.Range("E38").Value = "  -7.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.73"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -6.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "308.17"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -10.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.19"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -5.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.826"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -12.12%  "
$ws.Range("E43").Value = "  -8.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.997"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("E45").Value = "  -2.48%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.579"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -6.24%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0930"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0520"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -6.57%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "121.34"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -5.80%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0227"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -5.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.20"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -9.81%  "
